$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: localization handback succeeded, refresh status + widen columns to fit new text ---
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = 29.9777047293527
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527

# --- zh-cn sheet: status refreshed, handback datetime updated, stale error cleared ---
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("K2").Value = "2016-08-31 06:55:39"
# Assigning a plain "" clears the cell entirely; use a force-text empty string
# (leading apostrophe) then strip the resulting quote-prefix style so the
# cell survives as an empty shared-string value, same as the other empty
# cells on this row (e.g. L2/N2).
$zhcn.Range("P2").Value = "'"
$zhcn.Range("P2").Style = "Normal"
$zhcn.Columns.Item(3).ColumnWidth = 29.9777047293527
$zhcn.Columns.Item(16).ColumnWidth = 13.7470528738839

# --- de-de sheet: status refreshed, handback datetime updated, stale error cleared ---
$dede.Range("C2").Value = $newStatus
$dede.Range("K2").Value = "2016-08-31 06:55:46"
$dede.Range("P2").Value = "'"
$dede.Range("P2").Style = "Normal"
$dede.Columns.Item(3).ColumnWidth = 29.9777047293527
$dede.Columns.Item(16).ColumnWidth = 13.7470528738839
